# Update "想去人数" (F) and "最低票价" (G) figures on both the "展览" and
# "全部类型" worksheets, which contain identical data tables.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 347
    $ws.Range("F4").Value = 10614
    $ws.Range("F6").Value = 963
    $ws.Range("F7").Value = 110
    $ws.Range("F8").Value = 1311
    $ws.Range("F9").Value = 8224
    $ws.Range("G9").Value = 65
    $ws.Range("F10").Value = 30
    $ws.Range("F11").Value = 463
    $ws.Range("F15").Value = 3260
    $ws.Range("F17").Value = 326
    $ws.Range("F18").Value = 748
    $ws.Range("F20").Value = 1053
    $ws.Range("F23").Value = 1709
}
